$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 714.8276
$ws.Range("J17").Value = 803.12
$ws.Range("L17").Value = 2409.36
$ws.Range("N17").Value = -2745.36
$ws.Range("H51").Value = 51926.816
$ws.Range("J51").Value = 42119.5
$ws.Range("L51").Value = 42119.5
$ws.Range("N51").Value = -43087.5
$ws.Range("H70").Value = 8933.333000000001
$ws.Range("J70").Value = 8933.333000000001
$ws.Range("L70").Value = 26799.999
$ws.Range("N70").Value = -27339.999
$ws.Range("H73").Value = 8933.333000000001
$ws.Range("J73").Value = 8933.333000000001
$ws.Range("L73").Value = 26799.999
$ws.Range("N73").Value = -28671.999
$ws.Range("H98").Value = 125002376
$ws.Range("I98").Value = 142859460
$ws.Range("K98").Value = 142859460
$ws.Range("M98").Value = -142857962
$ws.Range("H122").Value = 125002376
$ws.Range("I122").Value = 142859460
$ws.Range("K122").Value = 428578380
$ws.Range("M122").Value = -428575930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H32").Value = 8773584
$ws.Range("I32").Value = 8773584
$ws.Range("K32").Value = 8773584
$ws.Range("M32").Value = -8773297
$ws.Range("H45").Value = 2319.3157
$ws.Range("I45").Value = 2624.6667
$ws.Range("J45").Value = 2262.0625
$ws.Range("K45").Value = 2624.6667
$ws.Range("L45").Value = 2262.0625
$ws.Range("M45").Value = -2247.6667
$ws.Range("N45").Value = -3016.0625
$ws.Range("H61").Value = 83509040
$ws.Range("I61").Value = 500000000
$ws.Range("J61").Value = 210849
$ws.Range("K61").Value = 500000000
$ws.Range("L61").Value = 210849
$ws.Range("M61").Value = -499999788
$ws.Range("N61").Value = -211273
$ws.Range("H74").Value = 7583287
$ws.Range("I74").Value = 11906619
$ws.Range("K74").Value = 11906619
$ws.Range("M74").Value = -11905745
$ws.Range("H77").Value = 7583287
$ws.Range("I77").Value = 11906619
$ws.Range("K77").Value = 59533095
$ws.Range("M77").Value = -59528727
$ws.Range("H132").Value = 5098.457
$ws.Range("I132").Value = 2303.3462
$ws.Range("J132").Value = 13173.223
$ws.Range("K132").Value = 6910.0386
$ws.Range("L132").Value = 39519.669
$ws.Range("M132").Value = -4380.0386
$ws.Range("N132").Value = -44579.669
$ws.Range("H136").Value = 83509040
$ws.Range("I136").Value = 500000000
$ws.Range("J136").Value = 210849
$ws.Range("K136").Value = 1500000000
$ws.Range("L136").Value = 632547
$ws.Range("M136").Value = -1499997450
$ws.Range("N136").Value = -637647

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H134").Value = 28184.756
$ws.Range("I134").Value = 3659.5264
$ws.Range("J134").Value = 338837.66
$ws.Range("K134").Value = 10978.5792
$ws.Range("L134").Value = 1016512.98
$ws.Range("M134").Value = -8443.5792
$ws.Range("N134").Value = -1021582.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1521.3914
$ws.Range("I7").Value = 113.76923
$ws.Range("K7").Value = 113.76923
$ws.Range("M7").Value = -0.7692299999999932
$ws.Range("H16").Value = 1888.3636
$ws.Range("I16").Value = 1553.2
$ws.Range("K16").Value = 1553.2
$ws.Range("M16").Value = -1266.2
$ws.Range("H31").Value = 1230149.4
$ws.Range("J31").Value = 1374577.4
$ws.Range("L31").Value = 1374577.4
$ws.Range("N31").Value = -1375167.4
$ws.Range("H34").Value = 1230149.4
$ws.Range("J34").Value = 1374577.4
$ws.Range("L34").Value = 1374577.4
$ws.Range("N34").Value = -1374981.4
$ws.Range("H55").Value = 25300
$ws.Range("I55").Value = 25300
$ws.Range("K55").Value = 25300
$ws.Range("M55").Value = -24985
$ws.Range("H113").Value = 1888.3636
$ws.Range("I113").Value = 1553.2
$ws.Range("K113").Value = 1553.2
$ws.Range("M113").Value = 616.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163.2
$ws.Range("J2").Value = 189.375
$ws.Range("L2").Value = 1136.25
$ws.Range("N2").Value = -1362.25
$ws.Range("H3").Value = 7960
$ws.Range("I3").Value = 1940
$ws.Range("K3").Value = 5820
$ws.Range("M3").Value = -5708
$ws.Range("H37").Value = 84000
$ws.Range("J37").Value = 84000
$ws.Range("L37").Value = 252000
$ws.Range("N37").Value = -252224
$ws.Range("H86").Value = 2648
$ws.Range("I86").Value = 728
$ws.Range("K86").Value = 2184
$ws.Range("M86").Value = -998
$ws.Range("H89").Value = 2648
$ws.Range("I89").Value = 728
$ws.Range("K89").Value = 6552
$ws.Range("M89").Value = -624
$ws.Range("H97").Value = 1367.6154
$ws.Range("I97").Value = 916.6667
$ws.Range("J97").Value = 1502.9
$ws.Range("K97").Value = 2750.0001
$ws.Range("L97").Value = 4508.700000000001
$ws.Range("M97").Value = -2254.0001
$ws.Range("N97").Value = -5500.700000000001
$ws.Range("H131").Value = 1323.5
$ws.Range("J131").Value = 1228.25
$ws.Range("L131").Value = 3684.75
$ws.Range("N131").Value = -13764.75
$ws.Range("H132").Value = 2409.7778
$ws.Range("J132").Value = 2833.8572
$ws.Range("L132").Value = 25504.7148
$ws.Range("N132").Value = -30564.7148
$ws.Range("H140").Value = 1422.4445
$ws.Range("I140").Value = 1422.4445
$ws.Range("K140").Value = 4267.333500000001
$ws.Range("M140").Value = 912.6664999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5169.5386
$ws.Range("I126").Value = 5315.7144
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 15947.1432
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -13477.1432
$ws.Range("N126").Value = -19937

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 59933.39
$ws.Range("I7").Value = 3458.8
$ws.Range("K7").Value = 3458.8
$ws.Range("M7").Value = -3346.8
$ws.Range("H22").Value = 1474.5
$ws.Range("I22").Value = 1474.5
$ws.Range("K22").Value = 1474.5
$ws.Range("M22").Value = -1179.5
$ws.Range("H27").Value = 1474.5
$ws.Range("I27").Value = 1474.5
$ws.Range("K27").Value = 1474.5
$ws.Range("M27").Value = -1367.5
$ws.Range("H100").Value = 3551.4285
$ws.Range("I100").Value = 3847.2727
$ws.Range("K100").Value = 3847.2727
$ws.Range("M100").Value = -3306.2727
$ws.Range("H126").Value = 59933.39
$ws.Range("I126").Value = 3458.8
$ws.Range("K126").Value = 10376.4
$ws.Range("M126").Value = -7906.400000000001
$ws.Range("H136").Value = 103072.87
$ws.Range("I136").Value = 128453.375
$ws.Range("K136").Value = 385360.125
$ws.Range("M136").Value = -382810.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22228778
$ws.Range("J62").Value = 25006500
$ws.Range("L62").Value = 25006500
$ws.Range("N62").Value = -25007748
$ws.Range("H65").Value = 22228778
$ws.Range("J65").Value = 25006500
$ws.Range("L65").Value = 125032500
$ws.Range("N65").Value = -125038740
$ws.Range("H81").Value = 100474
$ws.Range("I81").Value = 100474
$ws.Range("K81").Value = 200948
$ws.Range("M81").Value = -199887
$ws.Range("H84").Value = 100474
$ws.Range("I84").Value = 100474
$ws.Range("K84").Value = 1004740
$ws.Range("M84").Value = -999436
$ws.Range("H100").Value = 1859.5902
$ws.Range("J100").Value = 1980.7858
$ws.Range("L100").Value = 3961.5716
$ws.Range("N100").Value = -5043.5716
$ws.Range("H126").Value = 7840.933
$ws.Range("I126").Value = 7227.522
$ws.Range("K126").Value = 21682.566
$ws.Range("M126").Value = -19212.566
$ws.Range("H132").Value = 4574.7334
$ws.Range("I132").Value = 3875.3635
$ws.Range("J132").Value = 6498
$ws.Range("K132").Value = 11626.0905
$ws.Range("L132").Value = 19494
$ws.Range("M132").Value = -9096.0905
$ws.Range("N132").Value = -24554
$ws.Range("H136").Value = 1505
$ws.Range("I136").Value = 1505
$ws.Range("K136").Value = 4515
$ws.Range("M136").Value = -1965
